$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fndc5"
$ws.Range("C2").Value = "Itgav"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05565899999999999
$ws.Range("H2").Value = 0.166977
$ws.Range("I2").Value = 0.01431489546586528
$ws.Range("J2").Value = 0.01431489546586528
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 13.441269
$ws.Range("N2").Value = 40.323807
$ws.Range("O2").Value = 0.08973082133481231
$ws.Range("P2").Value = 0.08973082133481232
$ws.Range("Q2").Value = 0.7481275912709999
$ws.Range("R2").Value = 6.733148321439
$ws.Range("S2").Value = 0.001284487327474073
$ws.Range("T2").Value = 0.001284487327474073

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fndc5"
$ws.Range("C3").Value = "Itgav"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05565899999999999
$ws.Range("H3").Value = 0.166977
$ws.Range("I3").Value = 0.01431489546586528
$ws.Range("J3").Value = 0.01431489546586528
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 54.711535
$ws.Range("N3").Value = 164.134605
$ws.Range("O3").Value = 0.3652416280068742
$ws.Range("P3").Value = 0.3652416280068742
$ws.Range("Q3").Value = 3.045189326565
$ws.Range("R3").Value = 27.406703939085
$ws.Range("S3").Value = 0.005228395724700858
$ws.Range("T3").Value = 0.005228395724700859

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fndc5"
$ws.Range("C4").Value = "Itgav"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05565899999999999
$ws.Range("H4").Value = 0.166977
$ws.Range("I4").Value = 0.01431489546586528
$ws.Range("J4").Value = 0.01431489546586528
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 63.67711
$ws.Range("N4").Value = 191.03133
$ws.Range("O4").Value = 0.4250937452800914
$ws.Range("P4").Value = 0.4250937452800915
$ws.Range("Q4").Value = 3.544204265489999
$ws.Range("R4").Value = 31.89783838941
$ws.Range("S4").Value = 0.006085172526877673
$ws.Range("T4").Value = 0.006085172526877674

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Fndc5"
$ws.Range("C5").Value = "Itgav"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.05565899999999999
$ws.Range("H5").Value = 0.166977
$ws.Range("I5").Value = 0.01431489546586528
$ws.Range("J5").Value = 0.01431489546586528
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 17.96553866666667
$ws.Range("N5").Value = 53.896616
$ws.Range("O5").Value = 0.119933805378222
$ws.Range("P5").Value = 0.119933805378222
$ws.Range("Q5").Value = 0.9999439166479999
$ws.Range("R5").Value = 8.999495249832
$ws.Range("S5").Value = 0.00171683988681268
$ws.Range("T5").Value = 0.00171683988681268

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fndc5"
$ws.Range("C6").Value = "Itgav"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.907196
$ws.Range("H6").Value = 5.721588000000001
$ws.Range("I6").Value = 0.4905102745812252
$ws.Range("J6").Value = 0.4905102745812252
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 13.441269
$ws.Range("N6").Value = 40.323807
$ws.Range("O6").Value = 0.08973082133481231
$ws.Range("P6").Value = 0.08973082133481232
$ws.Range("Q6").Value = 25.635134471724
$ws.Range("R6").Value = 230.716210245516
$ws.Range("S6").Value = 0.04401388981133764
$ws.Range("T6").Value = 0.04401388981133765

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fndc5"
$ws.Range("C7").Value = "Itgav"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.907196
$ws.Range("H7").Value = 5.721588000000001
$ws.Range("I7").Value = 0.4905102745812252
$ws.Range("J7").Value = 0.4905102745812252
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 54.711535
$ws.Range("N7").Value = 164.134605
$ws.Range("O7").Value = 0.3652416280068742
$ws.Range("P7").Value = 0.3652416280068742
$ws.Range("Q7").Value = 104.34562070586
$ws.Range("R7").Value = 939.11058635274
$ws.Range("S7").Value = 0.1791547712421455
$ws.Range("T7").Value = 0.1791547712421456

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Fndc5"
$ws.Range("C8").Value = "Itgav"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.907196
$ws.Range("H8").Value = 5.721588000000001
$ws.Range("I8").Value = 0.4905102745812252
$ws.Range("J8").Value = 0.4905102745812252
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 63.67711
$ws.Range("N8").Value = 191.03133
$ws.Range("O8").Value = 0.4250937452800914
$ws.Range("P8").Value = 0.4250937452800915
$ws.Range("Q8").Value = 121.44472948356
$ws.Range("R8").Value = 1093.00256535204
$ws.Range("S8").Value = 0.208512849720099
$ws.Range("T8").Value = 0.2085128497200991

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Fndc5"
$ws.Range("C9").Value = "Itgav"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.907196
$ws.Range("H9").Value = 5.721588000000001
$ws.Range("I9").Value = 0.4905102745812252
$ws.Range("J9").Value = 0.4905102745812252
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 17.96553866666667
$ws.Range("N9").Value = 53.896616
$ws.Range("O9").Value = 0.119933805378222
$ws.Range("P9").Value = 0.119933805378222
$ws.Range("Q9").Value = 34.263803482912
$ws.Range("R9").Value = 308.374231346208
$ws.Range("S9").Value = 0.0588287638076429
$ws.Range("T9").Value = 0.0588287638076429

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Fndc5"
$ws.Range("C10").Value = "Itgav"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.3444803333333333
$ws.Range("H10").Value = 1.033441
$ws.Range("I10").Value = 0.08859663238134165
$ws.Range("J10").Value = 0.08859663238134165
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 13.441269
$ws.Range("N10").Value = 40.323807
$ws.Range("O10").Value = 0.08973082133481231
$ws.Range("P10").Value = 0.08973082133481232
$ws.Range("Q10").Value = 4.630252825542999
$ws.Range("R10").Value = 41.67227542988699
$ws.Range("S10").Value = 0.007949848591076213
$ws.Range("T10").Value = 0.007949848591076215

$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Fndc5"
$ws.Range("C11").Value = "Itgav"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.3444803333333333
$ws.Range("H11").Value = 1.033441
$ws.Range("I11").Value = 0.08859663238134165
$ws.Range("J11").Value = 0.08859663238134165
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 54.711535
$ws.Range("N11").Value = 164.134605
$ws.Range("O11").Value = 0.3652416280068742
$ws.Range("P11").Value = 0.3652416280068742
$ws.Range("Q11").Value = 18.84704781397833
$ws.Range("R11").Value = 169.623430325805
$ws.Range("S11").Value = 0.03235917824688777
$ws.Range("T11").Value = 0.03235917824688777

$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Fndc5"
$ws.Range("C12").Value = "Itgav"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.3444803333333333
$ws.Range("H12").Value = 1.033441
$ws.Range("I12").Value = 0.08859663238134165
$ws.Range("J12").Value = 0.08859663238134165
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 63.67711
$ws.Range("N12").Value = 191.03133
$ws.Range("O12").Value = 0.4250937452800914
$ws.Range("P12").Value = 0.4250937452800915
$ws.Range("Q12").Value = 21.93551207850333
$ws.Range("R12").Value = 197.41960870653
$ws.Range("S12").Value = 0.03766187427818794
$ws.Range("T12").Value = 0.03766187427818795

$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Fndc5"
$ws.Range("C13").Value = "Itgav"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.3444803333333333
$ws.Range("H13").Value = 1.033441
$ws.Range("I13").Value = 0.08859663238134165
$ws.Range("J13").Value = 0.08859663238134165
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 17.96553866666667
$ws.Range("N13").Value = 53.896616
$ws.Range("O13").Value = 0.119933805378222
$ws.Range("P13").Value = 0.119933805378222
$ws.Range("Q13").Value = 6.188774748406221
$ws.Range("R13").Value = 55.69897273565599
$ws.Range("S13").Value = 0.01062573126518971
$ws.Range("T13").Value = 0.01062573126518971

$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Fndc5"
$ws.Range("C14").Value = "Itgav"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.580852333333333
$ws.Range("H14").Value = 4.742557
$ws.Range("I14").Value = 0.4065781975715678
$ws.Range("J14").Value = 0.4065781975715678
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 13.441269
$ws.Range("N14").Value = 40.323807
$ws.Range("O14").Value = 0.08973082133481231
$ws.Range("P14").Value = 0.08973082133481232
$ws.Range("Q14").Value = 21.248661461611
$ws.Range("R14").Value = 191.237953154499
$ws.Range("S14").Value = 0.03648259560492437
$ws.Range("T14").Value = 0.03648259560492437

$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Fndc5"
$ws.Range("C15").Value = "Itgav"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.580852333333333
$ws.Range("H15").Value = 4.742557
$ws.Range("I15").Value = 0.4065781975715678
$ws.Range("J15").Value = 0.4065781975715678
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 54.711535
$ws.Range("N15").Value = 164.134605
$ws.Range("O15").Value = 0.3652416280068742
$ws.Range("P15").Value = 0.3652416280068742
$ws.Range("Q15").Value = 86.49085776499832
$ws.Range("R15").Value = 778.4177198849849
$ws.Range("S15").Value = 0.1484992827931399
$ws.Range("T15").Value = 0.14849928279314

$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Fndc5"
$ws.Range("C16").Value = "Itgav"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.580852333333333
$ws.Range("H16").Value = 4.742557
$ws.Range("I16").Value = 0.4065781975715678
$ws.Range("J16").Value = 0.4065781975715678
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 63.67711
$ws.Range("N16").Value = 191.03133
$ws.Range("O16").Value = 0.4250937452800914
$ws.Range("P16").Value = 0.4250937452800915
$ws.Range("Q16").Value = 100.6641079234233
$ws.Range("R16").Value = 905.97697131081
$ws.Range("S16").Value = 0.1728338487549267
$ws.Range("T16").Value = 0.1728338487549267

$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Fndc5"
$ws.Range("C17").Value = "Itgav"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1.580852333333333
$ws.Range("H17").Value = 4.742557
$ws.Range("I17").Value = 0.4065781975715678
$ws.Range("J17").Value = 0.4065781975715678
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 17.96553866666667
$ws.Range("N17").Value = 53.896616
$ws.Range("O17").Value = 0.119933805378222
$ws.Range("P17").Value = 0.119933805378222
$ws.Range("Q17").Value = 28.40086372079022
$ws.Range("R17").Value = 255.607773487112
$ws.Range("S17").Value = 0.0487624704185767
$ws.Range("T17").Value = 0.04876247041857671

